$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.258.08"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "'1.970.99"
$ws.Range("E3").Value = "  -6.23%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +1.07%  "
$ws.Range("D5").Value = "'329.33"
$ws.Range("E5").Value = "  -4.39%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").Value = "'0.5012"
$ws.Range("E7").Value = "  -5.22%  "
$ws.Range("D8").Value = "'0.4237"
$ws.Range("E8").Value = "  -4.40%  "
$ws.Range("D9").Value = "'54.38"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'0.09071"
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("D11").Value = "'1.105"
$ws.Range("E11").Value = "  -5.71%  "
$ws.Range("D12").Value = "'23.22"
$ws.Range("E12").Value = "  -6.35%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'7.925"
$ws.Range("E13").Value = "  -7.46%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.949.01"
$ws.Range("E14").Value = "  -10.68%  "
$ws.Range("D15").Value = "'6.448"
$ws.Range("E15").Value = "  -6.60%  "
$ws.Range("D16").Value = "'1.014"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "'0.00001107"
$ws.Range("E17").Value = "  -4.45%  "
$ws.Range("D18").Value = "'91.82"
$ws.Range("E18").Value = "  -9.77%  "
$ws.Range("D19").Value = "'0.06718"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'19.33"
$ws.Range("E20").Value = "  -8.71%  "
$ws.Range("D21").Value = "'1.012"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "'5.966"
$ws.Range("E22").Value = "  -5.65%  "
$ws.Range("D23").Value = "'29.267.28"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("D24").Value = "'11.94"
$ws.Range("E24").Value = "  -4.60%  "
$ws.Range("D25").Value = "'2.307"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'20.73"
$ws.Range("E26").Value = "  -5.09%  "
$ws.Range("D27").Value = "'156.33"
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").Value = "'6.287"
$ws.Range("E28").Value = "  -8.89%  "
$ws.Range("D29").Value = "'2.288"
$ws.Range("E29").Value = "  -8.89%  "
$ws.Range("D30").Value = "'127.28"
$ws.Range("E30").Value = "  -5.27%  "
$ws.Range("D31").Value = "'1.058"
$ws.Range("E31").Value = "  -6.98%  "
$ws.Range("D32").Value = "'0.09934"
$ws.Range("E32").Value = "  -5.87%  "
$ws.Range("D33").Value = "'1.535"
$ws.Range("E33").Value = "  -7.70%  "
$ws.Range("D34").Value = "'5.825"
$ws.Range("E34").Value = "  -6.86%  "
$ws.Range("D35").Value = "'3.747"
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("D36").Value = "'0.02436"
$ws.Range("E36").Value = "  -7.48%  "
$ws.Range("D37").Value = "'9.124"
$ws.Range("E37").Value = "  -10.04%  "
$ws.Range("D38").Value = "'0.06395"
$ws.Range("E38").Value = "  -6.14%  "
$ws.Range("D39").Value = "'1.292"
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("D40").Value = "'0.6513"
$ws.Range("E40").Value = "  -7.18%  "
$ws.Range("D41").Value = "'11.59"
$ws.Range("E41").Value = "  -9.12%  "
$ws.Range("E42").Value = "  -8.87%  "
$ws.Range("D43").Value = "'1.012"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "'0.6282"
$ws.Range("E44").Value = "  -7.96%  "
$ws.Range("D45").Value = "'13.50"
$ws.Range("E45").Value = "  -6.22%  "
$ws.Range("D46").Value = "'2.191"
$ws.Range("E46").Value = "  -5.88%  "
$ws.Range("D47").Value = "'1.293"
$ws.Range("E47").Value = "  -6.48%  "
$ws.Range("D48").Value = "'3.490"
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("D49").Value = "'0.00000000333"
$ws.Range("E49").Value = "  -6.10%  "
$ws.Range("D50").Value = "'0.06908"
$ws.Range("E50").Value = "  -6.09%  "
$ws.Range("D51").Value = "'1.121"
$ws.Range("E51").Value = "  -7.93%  "
